# Add "Gamma" (R) and "Weight decay" (S) columns, populate blanks for the
# existing 21 data rows, append 3 new training-run rows, and let the used
# range / dimension grow to A1:S25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. New header cells -------------------------------------------------
$ws.Cells.Item(1, 18).Value = "Gamma"
$ws.Cells.Item(1, 19).Value = "Weight decay"

# Match the bold/centered/bordered header style used by the rest of row 1
# (copy format only, so we reuse the existing style instead of registering
# a new one).
$ws.Range("Q1").Copy()
$ws.Range("R1:S1").PasteSpecial(-4122)

# ---- 2. Blank-but-present R/S cells for the 21 existing data rows -------
# Writing a lone apostrophe forces a (empty) text cell into existence
# instead of being dropped as a no-op, matching the empty placeholder cells
# Excel already uses in column Q (e.g. Q2). Re-pasting plain formatting
# from a bare data cell afterwards clears the quote-prefix flag so no new
# "show as text" marker is left behind.
$ws.Range("R2:S22").Value = "'"
$ws.Range("A2").Copy()
$ws.Range("R2:S22").PasteSpecial(-4122)

# ---- 3. Three new training-run rows --------------------------------------
$newRows = @(
    @("2024-1-4 16:44:33", 20, 64, 0.001, "SGD", "CEL", 12.8, 32, 1.1581, 1.1581, 56.8985, "FER2013", "cuda:0", 4, 2, "Stationær", 256.1767744999961, 0, 0),
    @("2024-1-4 16:53:8",  20, 64, 0.001, "SGD", "CEL", 12.5, 32, 1.2979, 1.1205, 58.5426, "FER2013", "cuda:0", 4, 2, "Stationær", 249.4619421000061, 0, 0),
    @("2024-1-4 17:42:41", 20, 64, 0.001, "SGD", "CEL", 12.6, 32, 1.2513, 1.1665, 57.5952, "FER2013", "cuda:0", 4, 2, "Stationær", 251.8056318999988, 0, 0)
)

$r = 23
foreach ($row in $newRows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}
